# RHEM_template.xlsx edit: add SAR (sodium adsorption ratio) input column
# and Avg SY / TDS output columns to support saline-scenario runs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inputs")

# ---------------------------------------------------------------------
# 1) Insert a new column at G for the "SAR" input (shifts old G:Q -> H:R,
#    carrying their values/styles with them, matching a native Excel
#    "Insert Sheet Columns" at column G).
# ---------------------------------------------------------------------
$ws.Columns("G:G").Insert()

# New SAR header (style matches the other Inputs headers, e.g. F1)
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "SAR"

# New SAR example value (style matches the neighboring data cell, e.g. F2)
$ws.Range("F2").Copy()
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("G2").Value = 0

# Give the new column's blank template row (row 3) the same style as its
# neighbors.
$ws.Range("F3").Copy()
$ws.Range("G3").PasteSpecial(-4122)

# The "Slope Steepness" example value moved from the old I2 (15.80000019073486)
# into J2 as part of the shift, but was also simplified to a plain 15.
$ws.Range("J2").Value = 15

# ---------------------------------------------------------------------
# 2) Relabel "Units" -> "Units (metric only)" in C1 (style changes to
#    match the other plain headers, e.g. F1).
# ---------------------------------------------------------------------
$ws.Range("F1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value = "Units (metric only)"

# ---------------------------------------------------------------------
# 3) "Biological Crusts Cover" moves from column S into the previously
#    unused column R (directly after "Litter Cover" in Q).
# ---------------------------------------------------------------------
$ws.Range("Q1").Copy()
$ws.Range("R1").PasteSpecial(-4122)
$ws.Range("R1").Value = "Biological Crusts Cover            ( % )"

$ws.Range("Q2").Copy()
$ws.Range("R2").PasteSpecial(-4122)
$ws.Range("R2").Value = 0

$ws.Range("Q3").Copy()
$ws.Range("R3").PasteSpecial(-4122)

# Clear the old column-S header/value now that it has moved to R.
$ws.Range("S1").Clear()
$ws.Range("S2").Clear()

# ---------------------------------------------------------------------
# 4) Add two new output columns after "Avg Soil Loss" (V): "Avg SY" (W)
#    and "TDS" (X) -- both needed for saline-scenario output reporting.
# ---------------------------------------------------------------------
$ws.Range("T1").Copy()
$ws.Range("W1").PasteSpecial(-4122)
$ws.Range("W1").Value = "Avg SY (ton/ha/year)"

$ws.Range("T1").Copy()
$ws.Range("X1").PasteSpecial(-4122)
$ws.Range("X1").Value = "TDS (ton/ha/year)"

# The column insert pushed the "NOTE" cell from Y1 to Z1; move it back to
# Y1 (directly after the new TDS column) and clear the vacated Z1.
$noteText = $ws.Range("Z1").Value()
$ws.Range("Z1").Copy()
$ws.Range("Y1").PasteSpecial(-4122)
$ws.Range("Y1").Value = $noteText
$ws.Range("Z1").Clear()

# ---------------------------------------------------------------------
# 5) Row 2's "Scenario Name"/"Scenario Description" example values are now
#    literally the column labels themselves, and the Units example (C2)
#    switches from "2 (English)" to "1 (Metric)".
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "Scenario Name"
$ws.Range("B2").Value = "Scenario Description"
$ws.Range("C2").Value = 1

# ---------------------------------------------------------------------
# 6) Update the "Uniform" defined name to track the Slope Shape example
#    cell, which shifted from H2 to I2 because of the new SAR column.
# ---------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Uniform") {
        $n.RefersTo = "=Inputs!`$I`$2"
    }
}

# ---------------------------------------------------------------------
# 7) Update the sheet selection to reflect where the new output columns
#    landed (matches the saved workbook view after the edit).
# ---------------------------------------------------------------------
$ws.Range("T2:X3").Select()
